$d = $word.ActiveDocument

# The last paragraph in the body is empty (just the paragraph mark); the
# diff adds a new run of text to it, carrying an eastAsia font hint plus
# explicit sz/szCs of 32 half-points (16pt).
$paragraphs = $d.Paragraphs
$lastParagraph = $paragraphs.Item($paragraphs.Count)

# Collapse a range to the very start of that (empty) paragraph, i.e.
# *before* its paragraph mark, so inserting content there adds a run to
# the existing paragraph instead of splitting off a brand-new one.
$insertionPoint = $d.Range($lastParagraph.Range.Start, $lastParagraph.Range.Start)

$xml = @'
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r>
              <w:rPr>
                <w:rFonts w:hint="eastAsia"/>
                <w:sz w:val="32"/>
                <w:szCs w:val="32"/>
              </w:rPr>
              <w:t>使用git创建分支简单又快捷</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$insertionPoint.InsertXML($xml)
